$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 81.95162343978882
$ws.Range("B2").Value = 50.19635820388794
$ws.Range("C2").Value = 23.62021422386169
$ws.Range("D2").Value = 24.95228934288025

$ws.Range("A3").Value = 48.13237404823303
$ws.Range("B3").Value = 63.16478085517883
$ws.Range("C3").Value = 32.04254674911499
$ws.Range("D3").Value = 274.609308719635
